$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the diff (price/volume refresh + Filecoin/Hedera row swap)
$ws.Range("D2").Value = '35.678.90'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '1.897.68'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.693'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.18'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '57.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.69%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0986'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.799'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.60%  '
$ws.Range("D15").Value = '2.174.92'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.05'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.72%  '
$ws.Range("D17").Value = '1.897.59'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '35.675.43'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '246.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  +4.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.02%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("E28").Value = '  +2.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0609'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.28%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +15.01%  '
$ws.Range("E36").Value = '  -16.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.857'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("E39").Value = '  +7.76%  '
$ws.Range("E40").Value = '  +7.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.52%  '
$ws.Range("D45").Value = '1.320.53'
$ws.Range("E45").Value = '  +1.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0813'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.81%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '42.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.45%  '
